$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the N2O row entirely; rows below shift up by one.
$ws.Rows.Item(4).Delete()

# The CH4 row (now row 3) becomes the new "Bio" row with its own data.
$ws.Range("A3").Value = "Bio"
$ws.Range("B3").Value = 262.1093884855168
$ws.Range("C3").Value = 279.1467580569499
$ws.Range("D3").Value = 431.1378158502502
$ws.Range("E3").Value = 325.0836925024324
$ws.Range("F3").Value = 123.783734336585
$ws.Range("G3").Value = 256.3340692038661
